$d = $word.ActiveDocument
$d.Content.Find.Execute("asserts", $false, $false, $false, $false, $false, $true, 1, $false, "assets", 2)
